$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new date column at DS (14-nov) ---
$ws = $wb.Worksheets.Item("Prix Spot")
$ws.Range("DS1").EntireColumn.Insert()
$ws.Range("DS1").Value = "14-nov"
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 123).Value = "-"
}

# --- Sheet "Gaz": append row 151 ---
# (force the date-looking string to stay plain text, matching the
#  existing column-A cells, instead of Excel's date auto-detection)
$gaz = $wb.Worksheets.Item("Gaz")
$gazDate = $gaz.Cells.Item(151, 1)
$gazDate.NumberFormat = "@"
$gazDate.Value = "2025-11-12"
$gazDate.ClearFormats()
$gaz.Cells.Item(151, 2).Value = 29

# --- Sheet "CO2": append row 151 ---
$co2 = $wb.Worksheets.Item("CO2")
$co2Date = $co2.Cells.Item(151, 1)
$co2Date.NumberFormat = "@"
$co2Date.Value = "2025-11-12"
$co2Date.ClearFormats()
$co2.Cells.Item(151, 2).Value = 81.75
